$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: replace values
$ws.Range("A2").Value = "JPY=X"
$ws.Range("B2").Value = "Macro"
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = "AV"
$ws.Range("G2").Value = "1d"

# Rows 3 and 4: clear contents, keep formatting/styles
$ws.Range("A3:G4").ClearContents()

# Update selection to D2
$ws.Range("D2").Select()
